# Jogos_do_Dia_Betfair_Back_Lay_2025-10-09.xlsx - daily refresh
#
# The sheet lists today's matches (one row per game, columns A:AO). This
# update replaces the whole fixture list with the latest Betfair odds pull:
# two matches were added (an early Swiss friendly at the top, a late
# Honduras-Costa Rica World-Cup qualifier at the bottom) and every other
# match's odds were refreshed, which also shuffles the row order because
# the sheet is kept sorted by kickoff Time.
#
# Simplest/most robust way to express that with COM automation is to just
# rewrite the whole A2:AO8 block in the final desired layout, rather than
# trying to replay inserts/moves - the net result is identical and there's
# no formula/formatting to preserve in the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds plain text dates like "2025-10-09". Excel's COM Value
# setter auto-parses strings that look like dates into date serials, so
# force the column to Text format first to keep them as literal strings
# (matches the other text columns, which are stored as plain text too).
$ws.Range("B2:B8").NumberFormat = "@"

# Row 2: Friendly Matches - Winterthur v FC Vaduz (new row)
$row2 = @('Friendly Matches','2025-10-09','08:30:00','Winterthur','FC Vaduz',1.02,1000,1.02,1000,1.02,1000,1.01,1.01,1.3,1.14,1.3,1.14,1.18,1.14,1.03,1.03,1.01,1.01,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000)

# Row 3: Brazilian Serie B - Ferroviaria v Chapecoense
$row3 = @('Brazilian Serie B','2025-10-09','19:00:00','Ferroviaria','Chapecoense',2.52,2.68,3,3.3,3.3,3.4,1.48,1.09,3,1.42,1.67,2.34,1.25,4.3,1.91,1.94,1.43,1.63,11.5,11,22,60,10.5,8,14,42,17,13,21,60,40,34,55,160,32,1000)

# Row 4: Brazilian Serie B - Coritiba v Atletico GO
$row4 = @('Brazilian Serie B','2025-10-09','19:30:00','Coritiba','Atletico GO',1.87,1.97,5.3,6,3.3,3.45,1.58,1.11,2.54,1.56,1.52,2.44,1.18,5.6,2.22,1.67,1.2,2.02,10.5,14,42,200,6.4,8,980,130,970,11.5,29,140,24,28,70,260,27,200)

# Row 5: Chilean Primera Division - Univ Catolica (Chile) v Nublense
$row5 = @('Chilean Primera Division','2025-10-09','20:30:00','Univ Catolica (Chile)','Nublense',1.9,1.99,4.4,5,3.4,4.1,1.42,1.09,3.1,1.41,1.72,2.2,1.27,4.1,2.02,1.86,1.25,2,980,980,980,130,980,980,980,70,980,980,980,100,980,980,980,180,980,100)

# Row 6: FIFA World Cup Qualifiers - Americas - Nicaragua v Haiti
$row6 = @('FIFA World Cup Qualifiers - Americas','2025-10-09','21:00:00','Nicaragua','Haiti',3.15,3.35,2.44,2.48,3.45,3.7,1.43,1.07,3.55,1.32,1.88,2.08,1.32,3.5,1.79,2.06,1.67,1.42,14.5,11,16,34,13,8.6,12,28,23,15,19.5,980,980,980,55,130,980,23)

# Row 7: Brazilian Serie B - Remo v Athletico-PR
$row7 = @('Brazilian Serie B','2025-10-09','21:35:00','Remo','Athletico-PR',3.65,4,2.18,2.38,3.1,3.4,1.56,1.11,2.8,1.48,1.6,2.4,1.22,4.8,2.06,1.82,1.72,1.33,9.6,8,13,1000,12,7.4,12,1000,1000,1000,1000,60,100,1000,1000,210,1000,1000)

# Row 8: FIFA World Cup Qualifiers - Americas - Honduras v Costa Rica (new row)
$row8 = @('FIFA World Cup Qualifiers - Americas','2025-10-09','23:00:00','Honduras','Costa Rica',2.8,980,1.01,980,1.09,980,1.46,1.06,1.11,1.36,1.24,1.36,1.18,1.05,1.03,1.03,1.01,1.38,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000)

$allRows = @($row2, $row3, $row4, $row5, $row6, $row7, $row8)

$numRows = $allRows.Count
$numCols = $row2.Count

$data = New-Object 'object[,]' $numRows,$numCols
for ($r = 0; $r -lt $numRows; $r++) {
    $srcRow = $allRows[$r]
    for ($c = 0; $c -lt $numCols; $c++) {
        $data[$r, $c] = $srcRow[$c]
    }
}

$ws.Range("A2:AO8").Value = $data

# The Text number-format on column B was only needed to stop the date
# strings from being auto-parsed while writing; drop it again afterwards
# so the cells are left in the workbook's normal (default) style, same as
# every other text cell on the sheet.
$ws.Range("B2:B8").ClearFormats()
